{"js": "// Replace the inline \"Ancillary Structures\" illustration image with a\n// plain hyperlink run whose visible text is the image's source URL.\n//\n// Before:\n//   <w:p pStyle=\"FirstParagraph\"><w:r><w:drawing>...</w:drawing></w:r></w:p>\n// After:\n//   <w:p pStyle=\"FirstParagraph\">\n//     <w:hyperlink r:id=\"...\">\n//       <w:r><w:rPr><w:rStyle w:val=\"Hyperlink\"/></w:rPr><w:t>https://...jpg...</w:t></w:r>\n//     </w:hyperlink>\n//   </w:p>\n\nconst url =\n  \"https://ura.gov.sg/-/media/Corporate/Guidelines/Development-control/Industrial/B202_Setbacks_Ancillary_Structures_Substation.jpg?h=100%25&w=100%25\";\n\nconst pictures = context.document.body.inlinePictures;\npictures.load(\"items\");\nawait context.sync();\n\nif (pictures.items.length > 0) {\n  // The picture lives alone in its own paragraph - grab that paragraph so\n  // we can drop the hyperlink text in exactly where the picture was.\n  const picture = pictures.items[0];\n  const paragraph = picture.paragraph;\n\n  // Remove the picture itself.\n  picture.delete();\n\n  // Insert the URL text in its place and turn it into a real hyperlink\n  // (Word automatically applies the \"Hyperlink\" character style).\n  const range = paragraph.insertText(url, Word.InsertLocation.replace);\n  range.hyperlink = url;\n\n  await context.sync();\n}\n", "ps1": "# Replace the inline \"Ancillary Structures\" illustration image with a\n# plain hyperlink whose visible text is the image's source URL.\n#\n# Before:\n#   <w:p pStyle=\"FirstParagraph\"><w:r><w:drawing>...</w:drawing></w:r></w:p>\n# After:\n#   <w:p pStyle=\"FirstParagraph\">\n#     <w:hyperlink r:id=\"...\">\n#       <w:r><w:rPr><w:rStyle w:val=\"Hyperlink\"/></w:rPr><w:t>https://...jpg...</w:t></w:r>\n#     </w:hyperlink>\n#   </w:p>\n\n$d = $word.ActiveDocument\n\n$url = \"https://ura.gov.sg/-/media/Corporate/Guidelines/Development-control/Industrial/B202_Setbacks_Ancillary_Structures_Substation.jpg?h=100%25&w=100%25\"\n\nif ($d.InlineShapes.Count -gt 0) {\n    $shape = $d.InlineShapes.Item(1)\n    $range = $shape.Range\n\n    # Remove the picture itself.\n    $shape.Delete()\n\n    # Insert the URL text where the picture used to be, then wrap it in a\n    # real hyperlink (Word applies the \"Hyperlink\" character style).\n    $range.Text = $url\n    $d.Hyperlinks.Add($range, $url) | Out-Null\n}\n"}
